$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")
$ws.Range("K2:L2").NumberFormat = "@"
$ws.Range("K2").Value = "1139am"
$ws.Range("L2").Value = "0209pm"
